# Update the NBA schedule table from the "Dec 26, 2022" slate to the
# "Dec 29, 2022" slate, and drop the last game row (table shrinks from
# 7 to 6 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 8) entirely, shifting the dimension
# from A1:O8 down to A1:O7.
$ws.Rows("8:8").Delete()

# Row 2 (game 0): Thunder @ Hornets, 7:00p, Spectrum Center
$ws.Range("C2").Value = 198
$ws.Range("D2").Value = "Thu, Dec 29, 2022"
$ws.Range("E2").Value = "7:00p"
$ws.Range("F2").Value = "Oklahoma City Thunder"
$ws.Range("H2").Value = "Charlotte Hornets"
$ws.Range("M2").Value = "Spectrum Center"
$ws.Range("O2").Value = "Oklahoma City Thunder"

# Row 3 (game 1): Cavaliers @ Pacers, 7:00p, Gainbridge Fieldhouse
$ws.Range("C3").Value = 199
$ws.Range("D3").Value = "Thu, Dec 29, 2022"
$ws.Range("E3").Value = "7:00p"
$ws.Range("F3").Value = "Cleveland Cavaliers"
$ws.Range("H3").Value = "Indiana Pacers"
$ws.Range("M3").Value = "Gainbridge Fieldhouse"
$ws.Range("O3").Value = "Indiana Pacers"

# Row 4 (game 2): Clippers @ Celtics, 7:30p, TD Garden
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = "Thu, Dec 29, 2022"
$ws.Range("E4").Value = "7:30p"
$ws.Range("F4").Value = "Los Angeles Clippers"
$ws.Range("H4").Value = "Boston Celtics"
$ws.Range("M4").Value = "TD Garden"
$ws.Range("O4").Value = "Los Angeles Clippers"

# Row 5 (game 3): Grizzlies @ Raptors, 7:30p, Scotiabank Arena
$ws.Range("C5").Value = 201
$ws.Range("D5").Value = "Thu, Dec 29, 2022"
$ws.Range("E5").Value = "7:30p"
$ws.Range("F5").Value = "Memphis Grizzlies"
$ws.Range("H5").Value = "Toronto Raptors"
$ws.Range("M5").Value = "Scotiabank Arena"
$ws.Range("O5").Value = "Toronto Raptors"

# Row 6 (game 4): Knicks @ Spurs, 8:00p, AT&T Center
$ws.Range("C6").Value = 202
$ws.Range("D6").Value = "Thu, Dec 29, 2022"
$ws.Range("E6").Value = "8:00p"
$ws.Range("F6").Value = "New York Knicks"
$ws.Range("H6").Value = "San Antonio Spurs"
$ws.Range("M6").Value = "AT&T Center"
$ws.Range("O6").Value = "San Antonio Spurs"

# Row 7 (game 5): Rockets @ Mavericks, 8:30p, American Airlines Center
$ws.Range("C7").Value = 203
$ws.Range("D7").Value = "Thu, Dec 29, 2022"
$ws.Range("E7").Value = "8:30p"
$ws.Range("F7").Value = "Houston Rockets"
$ws.Range("H7").Value = "Dallas Mavericks"
$ws.Range("M7").Value = "American Airlines Center"
$ws.Range("O7").Value = "Dallas Mavericks"
